# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns, and
# fix the TheSandbox / TrustWalletToken row ordering (rows 44-45).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "25.772.45"
    "E2"  = "  +0.28%  "
    "D3"  = "1.755.24"
    "E3"  = "  +1.22%  "
    "E4"  = "  +0.02%  "
    "D5"  = "237.46"
    "E5"  = "  -0.15%  "
    "E6"  = "  +0.00%  "
    "E7"  = "  +3.63%  "
    "D8"  = "40.57"
    "E8"  = "  -2.77%  "
    "D9"  = "0.2634"
    "E9"  = "  +9.22%  "
    "D10" = "0.06224"
    "E10" = "  +3.80%  "
    "D11" = "1.749.63"
    "E11" = "  +0.83%  "
    "E12" = "  +3.94%  "
    "D13" = "15.55"
    "E13" = "  +5.56%  "
    "D14" = "0.6056"
    "E14" = "  +3.69%  "
    "D15" = "78.38"
    "E15" = "  +2.18%  "
    "D16" = "4.451"
    "E16" = "  +1.72%  "
    "D17" = "1.001"
    "E17" = "  -0.04%  "
    "E18" = "  -0.07%  "
    "D19" = "25.821.07"
    "E19" = "  +0.35%  "
    "E20" = "  +3.44%  "
    "D21" = "0.000006787"
    "E21" = "  +7.46%  "
    "D22" = "1.975.94"
    "E22" = "  +0.78%  "
    "D23" = "4.060"
    "E23" = "  +4.05%  "
    "D24" = "8.192"
    "E24" = "  +4.45%  "
    "D25" = "5.182"
    "E25" = "  +1.48%  "
    "D26" = "137.97"
    "E26" = "  +1.73%  "
    "D27" = "1.472"
    "E27" = "  +2.38%  "
    "E28" = "  +5.85%  "
    "D29" = "1.812"
    "E29" = "  -1.70%  "
    "D30" = "102.64"
    "E30" = "  +2.38%  "
    "E31" = "  +1.35%  "
    "D32" = "3.704"
    "E32" = "  +2.09%  "
    "D33" = "3.398"
    "E33" = "  +3.73%  "
    "E34" = "  +1.28%  "
    "D35" = "0.9993"
    "E35" = "  -0.08%  "
    "D36" = "2.648"
    "E36" = "  -1.34%  "
    "D37" = "1.006"
    "E37" = "  -1.36%  "
    "D38" = "0.6023"
    "E38" = "  -0.54%  "
    "E39" = "  -2.25%  "
    "D40" = "1.967"
    "E40" = "  -4.58%  "
    "D41" = "0.01550"
    "E41" = "  +4.20%  "
    "D42" = "1.000"
    "E42" = "  -0.07%  "
    "D43" = "103.57"
    "E43" = "  +1.94%  "
    "B44" = "TrustWalletToken"
    "C44" = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
    "D44" = "0.7495"
    "E44" = "  -5.34%  "
    "B45" = "TheSandbox"
    "C45" = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
    "D45" = "0.3816"
    "E45" = "  +0.76%  "
    "D46" = "4.894"
    "E46" = "  -4.37%  "
    "D47" = "0.05484"
    "E47" = "  +7.96%  "
    "E48" = "  +4.89%  "
    "D49" = "5.969"
    "E49" = "  -1.49%  "
    "E50" = "  +1.83%  "
    "E51" = "  +0.71%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Leading apostrophe forces text (not numeric) interpretation for
    # values like "1.001" or "0.2634" that would otherwise be parsed as
    # numbers; ClearFormats drops the resulting quote-prefix style so the
    # cell's style index is left exactly as it was (style 0 / default).
    $cell.Value = "'" + $updates[$ref]
    $cell.ClearFormats()
}
